$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Centroid value corrections (re-measured positions) ---
$ws.Range("B16").Value = 170.625
$ws.Range("B17").Value = 170.67500000000001
$ws.Range("B18").Value = 132.97499999999999
$ws.Range("C18").Value = -73.25
$ws.Range("B25").Value = 170.42500000000001
$ws.Range("B26").Value = 170.46250000000001
$ws.Range("C26").Value = -77.275000000000006

# --- Apply center alignment across the whole sheet (matches the data-entry formatting pass) ---
$ws.Cells.Select()
$excel.Selection.HorizontalAlignment = -4108

# --- Restore the normal/default selection anchor ---
$ws.Range("F4").Select()
